$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on price (D-column) cells whose new value would otherwise
# be auto-coerced to a number by Excel, so they remain text like the original inline strings.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.299.23'
$ws.Range("E2").Value = '  +0.52%  '
$ws.Range("D3").Value = '1.601.87'
$ws.Range("E3").Value = '  +1.11%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '212.67'
$ws.Range("E5").Value = '  +0.68%  '
$ws.Range("D6").Value = '0.502'
$ws.Range("E6").Value = '  +0.13%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("E8").Value = '  +0.17%  '
$ws.Range("D10").Value = '18.93'
$ws.Range("E10").Value = '  -1.38%  '
$ws.Range("E11").Value = '  +0.90%  '
$ws.Range("D12").Value = '1.828.38'
$ws.Range("E12").Value = '  +1.18%  '
$ws.Range("D13").Value = '1.595.35'
$ws.Range("E13").Value = '  +0.56%  '
$ws.Range("E14").Value = '  +0.04%  '
$ws.Range("E15").Value = '  -2.07%  '
$ws.Range("D16").Value = '63.63'
$ws.Range("E16").Value = '  -0.63%  '
$ws.Range("D17").Value = '26.308.95'
$ws.Range("E17").Value = '  +0.52%  '
$ws.Range("D18").Value = '227.06'
$ws.Range("E18").Value = '  +6.56%  '
$ws.Range("E19").Value = '  -0.32%  '
$ws.Range("D20").Value = '7.64'
$ws.Range("E20").Value = '  +4.27%  '
$ws.Range("D22").Value = '4.31'
$ws.Range("E22").Value = '  +1.89%  '
$ws.Range("D23").Value = '2.17'
$ws.Range("E23").Value = '  +0.30%  '
$ws.Range("D24").Value = '8.96'
$ws.Range("E24").Value = '  +0.11%  '
$ws.Range("D25").Value = '145.48'
$ws.Range("E25").Value = '  +1.33%  '
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("E27").Value = '  -0.22%  '
$ws.Range("E28").Value = '  +1.55%  '
$ws.Range("E29").Value = '  +2.25%  '
$ws.Range("E30").Value = '  -0.39%  '
$ws.Range("E31").Value = '  +1.00%  '
$ws.Range("D32").Value = '3.21'
$ws.Range("E32").Value = '  +0.72%  '
$ws.Range("D33").Value = '1.442.00'
$ws.Range("E33").Value = '  +7.61%  '
$ws.Range("D34").Value = '2.96'
$ws.Range("E34").Value = '  +1.19%  '
$ws.Range("E35").Value = '  -0.44%  '
$ws.Range("D37").Value = '0.566'
$ws.Range("E37").Value = '  -2.39%  '
$ws.Range("E38").Value = '  -0.88%  '
$ws.Range("D39").Value = '0.824'
$ws.Range("E39").Value = '  +1.12%  '
$ws.Range("E40").Value = '  +0.73%  '
$ws.Range("E41").Value = '  +0.13%  '
$ws.Range("E42").Value = '  +2.24%  '
$ws.Range("E43").Value = '  -2.57%  '
$ws.Range("D44").Value = '1.740.22'
$ws.Range("E44").Value = '  +1.20%  '
$ws.Range("D45").Value = '0.758'
$ws.Range("E45").Value = '  -1.29%  '
$ws.Range("D46").Value = '60.84'
$ws.Range("E46").Value = '  -0.23%  '
$ws.Range("D47").Value = '87.68'
$ws.Range("E47").Value = '  +2.05%  '
$ws.Range("E48").Value = '  +0.63%  '
$ws.Range("D49").Value = '0.0499'
$ws.Range("E49").Value = '  -0.15%  '
$ws.Range("D50").Value = '0.0₇0963'
$ws.Range("E50").Value = '  -7.54%  '
$ws.Range("D51").Value = '0.0952'
$ws.Range("E51").Value = '  -3.20%  '
